$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force the cell to remain plain text (matching the source's inline-string
    # cells) even when the value looks numeric (e.g. "0.577", "27.95"),
    # without leaving a residual custom style on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" '64.405.40'
Set-TextValue "E2" '  +0.09%  '
Set-TextValue "D3" '3.421.71'
Set-TextValue "E3" '  -1.14%  '
Set-TextValue "E4" '  +0.01%  '
Set-TextValue "D5" '571.33'
Set-TextValue "E5" '  -2.01%  '
Set-TextValue "D6" '159.79'
Set-TextValue "E6" '  +1.02%  '
Set-TextValue "E7" '  +0.09%  '
Set-TextValue "D8" '3.425.64'
Set-TextValue "E8" '  -1.24%  '
Set-TextValue "D9" '0.577'
Set-TextValue "E9" '  +7.96%  '
Set-TextValue "D10" '7.28'
Set-TextValue "E10" '  -4.34%  '
Set-TextValue "E11" '  +0.22%  '
Set-TextValue "E12" '  -1.21%  '
Set-TextValue "D13" '4.014.21'
Set-TextValue "E13" '  -1.09%  '
Set-TextValue "E14" '  -2.22%  '
Set-TextValue "E15" '  +1.65%  '
Set-TextValue "D16" '27.95'
Set-TextValue "E16" '  +1.14%  '
Set-TextValue "D17" '64.458.99'
Set-TextValue "E17" '  +0.16%  '
Set-TextValue "D18" '3.422.99'
Set-TextValue "E18" '  -0.98%  '
Set-TextValue "E19" '  -2.38%  '
Set-TextValue "D20" '14.13'
Set-TextValue "E20" '  -2.21%  '
Set-TextValue "D21" '383.74'
Set-TextValue "E21" '  -3.59%  '
Set-TextValue "D22" '8.11'
Set-TextValue "E22" '  -5.13%  '
Set-TextValue "D23" '72.71'
Set-TextValue "E23" '  +0.64%  '
Set-TextValue "B24" 'Polygon'
Set-TextValue "C24" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D24" '0.541'
Set-TextValue "E24" '  -1.13%  '
Set-TextValue "B25" 'Dai'
Set-TextValue "C25" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D25" '0.998'
Set-TextValue "E25" '  -0.20%  '
Set-TextValue "E26" '  +9.74%  '
Set-TextValue "E27" '  +0.10%  '
Set-TextValue "E28" '  -1.42%  '
Set-TextValue "E29" '  +0.22%  '
Set-TextValue "E30" '  +3.72%  '
Set-TextValue "E31" '  -0.25%  '
Set-TextValue "E32" '  -1.25%  '
Set-TextValue "E33" '  -2.71%  '
Set-TextValue "D34" '23.46'
Set-TextValue "E34" '  -1.55%  '
Set-TextValue "E35" '  +0.09%  '
Set-TextValue "D36" '7.06'
Set-TextValue "E36" '  +2.55%  '
Set-TextValue "D37" '162.03'
Set-TextValue "E37" '  +2.10%  '
Set-TextValue "E38" '  -1.86%  '
Set-TextValue "D39" '3.014.44'
Set-TextValue "E39" '  +5.04%  '
Set-TextValue "D40" '1.89'
Set-TextValue "E40" '  +0.43%  '
Set-TextValue "D41" '0.0761'
Set-TextValue "E41" '  -3.25%  '
Set-TextValue "D42" '26.99'
Set-TextValue "E42" '  -5.35%  '
Set-TextValue "E43" '  +2.01%  '
Set-TextValue "E44" '  -2.29%  '
Set-TextValue "D45" '42.56'
Set-TextValue "E45" '  +0.84%  '
Set-TextValue "E46" '  -2.22%  '
Set-TextValue "D47" '24.34'
Set-TextValue "E47" '  +6.67%  '
Set-TextValue "D48" '1.08'
Set-TextValue "E48" '  -3.05%  '
Set-TextValue "D49" '0.868'
Set-TextValue "E49" '  +3.55%  '
Set-TextValue "D50" '6.59'
Set-TextValue "E50" '  +2.65%  '
Set-TextValue "D51" '2.15'
Set-TextValue "E51" '  +1.05%  '
